$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates (numeric, column A) and moods (column G, reuse existing "Tired" string)
$ws.Range("A23").Value = 20200213
$ws.Range("A24").Value = 20200218
$ws.Range("A25").Value = 20200219
$ws.Range("G23").Value = "Tired"
$ws.Range("G24").Value = "Tired"
$ws.Range("G25").Value = "Tired"
$ws.Range("C23").Value = "Andre, Kaj"
$ws.Range("B23").Value = "5-8pm"

# New text cells, entered in the order that reproduces the author's
# original shared-string allocation order.
$ws.Range("D23").Value = "1. Midterm;`n2. Learn more KEPs;`n3. Learn about Stakeholders."
$ws.Range("E23").Value = "3 KEPs learned:`n1. Prioritize Stakeholders;`n2. Move along levels of abstraction;`n3. Do something else."
$ws.Range("F23").Value = "The exam has a theoretical part and a pratical part. I think the theoretical part was not difficult if you have gone through all the slides and memorized some key concepts. However, the practical part of analyzing the pacman code seemed tricky because the question might have different interpretations. In the first question we were asked to use SimpleUML to draw a diagram of classes DIRECTLY related to the Game class and the SinglePlayerGame class without writing fields and constructors (what about methods?). Here the word DIRECTLY did not have a clear definition. In the next question we had to purpose some beacons, which also varied from person to person. So I was confused if I really understood what the questions wanted."

$ws.Range("B24").Value = "3:00-3:20pm"
$ws.Range("C24").Value = "Kaj, Wenchia"
$ws.Range("D24").Value = "To discuss what we can do to improve our report for homework 2"

$ws.Range("B25").Value = "10:00-11:00am"
$ws.Range("C25").Value = "Kaj, Wenchia"

$ws.Range("E24").Value = "We learned that the features picked by us in the report were not essential in the project. And we were asked to what was inside the realm-java and figure out was it really just a wrapper or something more."

$ws.Range("D25").Value = "To report our discovery to Kaj and see what we can do"

$ws.Range("F24").Value = "It was very frustrating to know why our homework got such a low score. We even thought about switching to different projects but knew that was not likely to happen. So instead we had to dig into realm-core written in C++, making the whole process more challenging. We might not be on the right track of searching essential features because they are not even implemented in realm-java."

$ws.Range("E25").Value = "We still cannot confirm whether realm-java is only a wrapper or something more. But we received suggestions from Kaj that we should investigate what realm-java provides differently from the other databases like SQLite. "

$ws.Range("F25").Value = "We thought we should know how to use realm-java in real projects first, then try to understand what makes realm-java stands out among all these database libraries. If it is just nothing new, then we will have to research essential features like querying in realm-core."

# Row heights to match the author's manual adjustments for the new rows
$ws.Rows.Item(23).RowHeight = 399.6
$ws.Rows.Item(24).RowHeight = 232.8
$ws.Rows.Item(25).RowHeight = 154.2

# Leave the cursor/selection where the author left it after typing the new entries
$ws.Range("G26").Select()
